$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (Method Inputs / Condition) rows 7-14
$ws.Cells.Item(7, 5).Value = "No account exists. Create account with valid data."
$null = $ws.Range("D7").Copy()
$null = $ws.Range("E7").PasteSpecial(-4122)
$ws.Cells.Item(8, 5).Value = "Account created with invalid overdraft limit."
$ws.Cells.Item(9, 5).Value = "Account created with invalid overdraft rate."
$ws.Cells.Item(10, 5).Value = "Account created with invalid date type."
$ws.Cells.Item(11, 5).Value = "Account with balance = 0, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Cells.Item(12, 5).Value = "Account with balance = -600, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Cells.Item(13, 5).Value = "Account with balance = -100, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Cells.Item(14, 5).Value = "Account with valid balance and overdraft details"

# Column F (Method Inputs) rows 7-14
$ws.Cells.Item(7, 6).Value = "account_number=1234567, client_number=22, balance=1000.00, date_created=today, overdraft_limit=-100.0, overdraft_rate=0.05"
$ws.Cells.Item(8, 6).Value = 'overdraft_limit="invalid"'
$ws.Cells.Item(9, 6).Value = 'overdraft_rate="invalid"'
$ws.Cells.Item(10, 6).Value = 'date_created="2025-10-27"'
$ws.Cells.Item(11, 6).Value = "Call get_service_charges()"
$ws.Cells.Item(12, 6).Value = "Call get_service_charges()"
$ws.Cells.Item(13, 6).Value = "Call get_service_charges()"
$ws.Cells.Item(14, 6).Value = "Call str(account)"

# Column G (Expected Result) rows 7, 8, 9, 11, 12
$ws.Cells.Item(7, 7).Value = "All attributes correctly set. Balance = 1000.00, overdraft_limit = -100.0, overdraft_rate = 0.05."
$ws.Cells.Item(8, 7).Value = "overdraft_limit defaults to -100.0"
$ws.Cells.Item(9, 7).Value = "overdraft_rate defaults to 0.05"
$ws.Cells.Item(11, 7).Value = "Expected = 0.50"
$ws.Cells.Item(12, 7).Value = "Expected = 25.50"

# Developer name
$ws.Range("C3").Value = "Parneet kaur"

# Remaining column G cells
$ws.Cells.Item(10, 7).Value = "date_created defaults to today's date (2025-10-27)"
$ws.Cells.Item(13, 7).Value = "Expected = 0.50"
$ws.Cells.Item(14, 7).Value = "Returns string: Account Number: 1234567 Balance: 1000.00 Overdraft Limit: -100.00 Overdraft Rate: 5.00% Account Type: Chequing"

# Update the view: scroll so row 2 is the top visible row, and select J13
$null = $ws.Activate()
$null = $ws.Range("J13").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
